$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.747.23'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '3.075.71'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').Value = '''239.55'
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').Value = '''617.02'
$ws.Range('E6').Value = '  -1.55%  '
$ws.Range('E7').Value = '  +1.50%  '
$ws.Range('D8').Value = '''0.362'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '3.074.25'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('E12').Value = '  +2.81%  '
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').Value = '''34.42'
$ws.Range('E14').Value = '  -5.34%  '
$ws.Range('D15').Value = '''5.44'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').Value = '89.697.17'
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('D17').Value = '3.652.06'
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').Value = '3.087.13'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('E19').Value = '  -4.00%  '
$ws.Range('E20').Value = '  +2.44%  '
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').Value = '''5.74'
$ws.Range('E22').Value = '  +3.08%  '
$ws.Range('D23').Value = '''435.12'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = '''5.59'
$ws.Range('E25').Value = '  -4.89%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '''90.26'
$ws.Range('E26').Value = '  +1.13%  '
$ws.Range('D27').Value = '''11.71'
$ws.Range('E27').Value = '  -4.77%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '''0.242'
$ws.Range('E30').Value = '  +18.49%  '
$ws.Range('D31').Value = '''0.176'
$ws.Range('E31').Value = '  +10.34%  '
$ws.Range('D32').Value = '''0.116'
$ws.Range('E32').Value = '  +28.14%  '
$ws.Range('D33').Value = '''9.07'
$ws.Range('E33').Value = '  -4.18%  '
$ws.Range('D34').Value = '''0.968'
$ws.Range('E34').Value = '  +8.19%  '
$ws.Range('E35').Value = '  +8.87%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').Value = '''7.64'
$ws.Range('E36').Value = '  +8.34%  '
$ws.Range('B37').Value = 'MantraDAO'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D37').Value = '''4.28'
$ws.Range('E37').Value = '  +23.92%  '
$ws.Range('D38').Value = '''26.06'
$ws.Range('E38').Value = '  -1.24%  '
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').Value = '''482.27'
$ws.Range('E40').Value = '  -5.35%  '
$ws.Range('D41').Value = '''3.49'
$ws.Range('E41').Value = '  -8.14%  '
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('D43').Value = '''0.414'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '''154.52'
$ws.Range('E46').Value = '  +2.18%  '
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('D48').Value = '''0.678'
$ws.Range('E48').Value = '  -1.41%  '
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('D50').Value = '''44.01'
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('D51').Value = '''0.998'
$ws.Range('E51').Value = '  -0.13%  '
